# Append three new list items (same ListParagraph / numId=3 bullet list)
# right after the existing "...How can Angular consume services?" item,
# at the very end of the document body.
#
# We use Find & Replace to grow the last list paragraph into five
# paragraphs (original text + the three new questions + a temporary
# sentinel paragraph). The sentinel is needed because Word only copies
# run-level character formatting (the themed fonts / size) from the
# paragraph mark that follows the insertion point; at the true end of
# the story there is no following paragraph mark to copy from, which
# would otherwise leave the very last new run without its w:rPr. Adding
# the sentinel paragraph gives the last real run something to inherit
# from, and we then delete the sentinel (together with the paragraph
# break that created it) so the document returns to its original
# paragraph-count shape plus exactly three new items.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "How can Angular consume services?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "How can Angular consume services?^p" + `
    "What are forms in Angular?^p" + `
    "How do forms consume information?^p" + `
    "What is Reactive programming?^p" + `
    "ZZZ_SENTINEL_ZZZ",
    2) | Out-Null

# Locate the sentinel paragraph and drop it, along with the paragraph
# mark that separates it from "What is Reactive programming?", so the
# last real paragraph again ends the document (picking up the themed
# run formatting along the way).
$sentinel = $d.Content
$sentinel.Find.Execute("ZZZ_SENTINEL_ZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sentinel.Delete()

$trailingMark = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$trailingMark.Delete()
